$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, date range) ---
$txt = $ws.Range("A8").Text
$idx = $txt.IndexOf("48")
$ws.Range("A8").Characters($idx + 1, 2).Text = "49"

$txt = $ws.Range("C9").Text
$idx = $txt.IndexOf("11/27/2023")
$ws.Range("C9").Characters($idx + 1, 10).Text = "12/4/2023"
$txt = $ws.Range("C9").Text
$idx = $txt.IndexOf("12/3/2023")
$ws.Range("C9").Characters($idx + 1, 9).Text = "12/10/2023"

# --- Cell type changes (number <-> shared string) via Copy to preserve type+style ---
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C16").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 2
$ws.Range("H14").Copy($ws.Range("E27"))

# --- Numeric value updates ---
$ws.Range("L15").Value = -41.025641025641
$ws.Range("N15").Value = -42.5
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 16.666666666666
$ws.Range("I16").Value = 163
$ws.Range("J16").Value = 171
$ws.Range("K16").Value = -4.678362573099
$ws.Range("L16").Value = 20.74074074074
$ws.Range("M16").Value = -48.902821316614
$ws.Range("N16").Value = -83.485309017223
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 3.703703703703
$ws.Range("I17").Value = 413
$ws.Range("J17").Value = 401
$ws.Range("K17").Value = 2.992518703241
$ws.Range("L17").Value = 7.272727272727
$ws.Range("M17").Value = 55.849056603773
$ws.Range("N17").Value = 5.626598465473
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -28.571428571428
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 5.263157894736
$ws.Range("I18").Value = 222
$ws.Range("J18").Value = 204
$ws.Range("K18").Value = 8.823529411764
$ws.Range("L18").Value = 45.098039215686
$ws.Range("M18").Value = -40.161725067385
$ws.Range("N18").Value = -85.73264781491
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = 7.54716981132
$ws.Range("I19").Value = 653
$ws.Range("J19").Value = 599
$ws.Range("K19").Value = 9.015025041736
$ws.Range("L19").Value = 41.648590021692
$ws.Range("M19").Value = 51.157407407407
$ws.Range("N19").Value = 16.399286987522
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -70
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 326
$ws.Range("J20").Value = 262
$ws.Range("K20").Value = 24.42748091603
$ws.Range("L20").Value = 117.333333333333
$ws.Range("M20").Value = -10.43956043956
$ws.Range("N20").Value = -90.10321797207
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = -26.190476190476
$ws.Range("F21").Value = 143
$ws.Range("G21").Value = 139
$ws.Range("H21").Value = 2.877697841726
$ws.Range("I21").Value = 1802
$ws.Range("J21").Value = 1660
$ws.Range("K21").Value = 8.554216867469
$ws.Range("L21").Value = 35.692771084337
$ws.Range("M21").Value = 1.008968609865
$ws.Range("N21").Value = -73.697270471464
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -26.470588235294
$ws.Range("F24").Value = 85
$ws.Range("H24").Value = -8.602150537634
$ws.Range("I24").Value = 1179
$ws.Range("J24").Value = 1428
$ws.Range("K24").Value = -17.436974789916
$ws.Range("L24").Value = 16.272189349112
$ws.Range("M24").Value = 53.715775749674
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 27.272727272727
$ws.Range("F25").Value = 52
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = 26.829268292682
$ws.Range("I25").Value = 659
$ws.Range("J25").Value = 616
$ws.Range("K25").Value = 6.980519480519
$ws.Range("L25").Value = 37.291666666666
$ws.Range("M25").Value = 15.20979020979
$ws.Range("C26").Value = 1
$ws.Range("I26").Value = 40
$ws.Range("K26").Value = 37.931034482758
$ws.Range("L26").Value = -23.076923076923
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -85.714285714285
$ws.Range("I27").Value = 44
$ws.Range("J27").Value = 59
$ws.Range("K27").Value = -25.423728813559
$ws.Range("L27").Value = -10.204081632653
$ws.Range("N28").Value = -89.743589743589
$ws.Range("N29").Value = -89.855072463768
